# Append the new data row (row 95) to Sheet1, mirroring the existing
# "date / weekday / hour / ranking" rows above it.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 95

# Leading apostrophe forces these to be stored as literal text instead of
# being auto-parsed into a date serial number (matches columns A/B of all
# the other rows, which are plain text, not real Excel dates).
$ws.Range("A$newRow").Value = "'2025/10/12"
$ws.Range("B$newRow").Value = "'日"

# Numeric columns.
$ws.Range("C$newRow").Value = 20
$ws.Range("D$newRow").Value = 201

# Writing literal text via the apostrophe prefix makes Excel tag the cell
# with a "quote prefix" style; reset back to the default/Normal style so
# the new row carries no formatting, same as the rows above it.
$ws.Range("A$newRow`:B$newRow").Style = "Normal"
